$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 0.371
$ws.Range("C4").Value = 2.897
$ws.Range("D4").Value = 7.328
$ws.Range("E4").Value = 2.22
$ws.Range("F4").Value = 3.802
$ws.Range("G4").Value = 4.189
$ws.Range("H4").Value = 4.546
$ws.Range("B5").Value = 46.455
$ws.Range("C5").Value = 109172085.958
$ws.Range("D5").Value = 274556070.182
$ws.Range("E5").Value = 115084333.978
$ws.Range("F5").Value = 32.878
$ws.Range("G5").Value = 35.162
$ws.Range("H5").Value = 36.928
$ws.Range("B6").Value = 423.671
$ws.Range("C6").Value = 300610421.643
$ws.Range("D6").Value = 488700241.637
$ws.Range("E6").Value = 201811437.079
$ws.Range("F6").Value = 90.601
$ws.Range("G6").Value = 92.941
$ws.Range("H6").Value = 97.233
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 1020.734
$ws.Range("D10").Value = 6461.495
$ws.Range("E10").Value = 2129.866
$ws.Range("F10").Value = 1.594
$ws.Range("G10").Value = 1.804
$ws.Range("H10").Value = 2.09
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 30204.419
$ws.Range("D11").Value = 69741.77800000001
$ws.Range("E11").Value = 30895.739
$ws.Range("F11").Value = 12.651
$ws.Range("G11").Value = 13.408
$ws.Range("H11").Value = 15.106
$ws.Range("B12").Value = 108094.892
$ws.Range("C12").Value = 116438.343
$ws.Range("D12").Value = 123925.061
$ws.Range("E12").Value = 5939.353
$ws.Range("F12").Value = 30.852
$ws.Range("G12").Value = 34.752
$ws.Range("H12").Value = 37.93
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 3.738
$ws.Range("D16").Value = 19.317
$ws.Range("E16").Value = 7.481
$ws.Range("F16").Value = 3.478
$ws.Range("G16").Value = 3.738
$ws.Range("H16").Value = 4.257
$ws.Range("B17").Value = 16.644
$ws.Range("C17").Value = 19.144
$ws.Range("D17").Value = 20.699
$ws.Range("E17").Value = 1.339
$ws.Range("F17").Value = 23.336
$ws.Range("G17").Value = 24.62
$ws.Range("H17").Value = 25.715
$ws.Range("B18").Value = 19.094
$ws.Range("C18").Value = 20.136
$ws.Range("D18").Value = 20.863
$ws.Range("E18").Value = 0.75
$ws.Range("F18").Value = 60.857
$ws.Range("G18").Value = 64.611
$ws.Range("H18").Value = 75.34099999999999
$ws.Range("B22").Value = 0.034
$ws.Range("C22").Value = 37.907
$ws.Range("D22").Value = 138.325
$ws.Range("E22").Value = 49.55
$ws.Range("F22").Value = 3.482
$ws.Range("G22").Value = 3.741
$ws.Range("H22").Value = 4.165
$ws.Range("B23").Value = 0.038
$ws.Range("C23").Value = 98.67700000000001
$ws.Range("D23").Value = 546.926
$ws.Range("E23").Value = 197.694
$ws.Range("F23").Value = 24.411
$ws.Range("G23").Value = 25.625
$ws.Range("H23").Value = 26.942
$ws.Range("B24").Value = 0.511
$ws.Range("C24").Value = 719.646
$ws.Range("D24").Value = 1123.076
$ws.Range("E24").Value = 476.124
$ws.Range("F24").Value = 68.542
$ws.Range("G24").Value = 75.956
$ws.Range("H24").Value = 87.59099999999999
$ws.Range("B28").Value = 10.945
$ws.Range("C28").Value = 21074.108
$ws.Range("D28").Value = 210531.148
$ws.Range("E28").Value = 63152.347
$ws.Range("F28").Value = 3.291
$ws.Range("G28").Value = 3.429
$ws.Range("H28").Value = 3.678
$ws.Range("B29").Value = 130.234
$ws.Range("C29").Value = 485757.737
$ws.Range("D29").Value = 1663701.346
$ws.Range("E29").Value = 741979.458
$ws.Range("F29").Value = 24.223
$ws.Range("G29").Value = 24.993
$ws.Range("H29").Value = 25.896
$ws.Range("B30").Value = 36074.024
$ws.Range("C30").Value = 2738934.314
$ws.Range("D30").Value = 3410228.831
$ws.Range("E30").Value = 930837.737
$ws.Range("F30").Value = 63.866
$ws.Range("G30").Value = 67.88500000000001
$ws.Range("H30").Value = 71.50700000000001
$ws.Range("B34").Value = -3.065
$ws.Range("C34").Value = 2752550.401
$ws.Range("D34").Value = 26766793.87
$ws.Range("E34").Value = 8007943.506
$ws.Range("F34").Value = 5.971
$ws.Range("G34").Value = 6.557
$ws.Range("H34").Value = 7.31
$ws.Range("B35").Value = -1.018
$ws.Range("C35").Value = 316285208.576
$ws.Range("D35").Value = 585142659.725
$ws.Range("E35").Value = 231936289.151
$ws.Range("F35").Value = 48.039
$ws.Range("G35").Value = 56.525
$ws.Range("H35").Value = 63.648
$ws.Range("B36").Value = 1.539
$ws.Range("C36").Value = 765617344.123
$ws.Range("D36").Value = 1395377722.226
$ws.Range("E36").Value = 527510245.444
$ws.Range("F36").Value = 129.332
$ws.Range("G36").Value = 145.234
$ws.Range("H36").Value = 165.223
